$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 333 ("nafo informado"), shifting rows below it up.
$ws.Rows.Item(333).Delete()
